$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "19/02/2020"
$ws.Range("B7").Value = "30k"
$ws.Range("A8").Value = "21/02/2020"
$ws.Range("B8").Value = "30k"

$ws.Columns.Item(2).ColumnWidth = 8.5

$ws.Range("B3").Select()
